$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize separators in research/position text cells: semicolons -> commas,
# full-width commas -> commas, and trim stray whitespace.
$ws.Range("I4").Value = "Netease Fuxi AI Lab, Alibaba"
$ws.Range("I5").Value = "Postdoctoral Fellow of UdeM/MILA, Associate Reseacher, Tianjin University"
$ws.Range("G6").Value = "Reinforcement Learning, Transfer Learning, Multiagent Learning"
$ws.Range("I6").Value = "Postdoc at University of Alberta, Associate Professor, Nanjing University"
$ws.Range("G7").Value = "Multiagent Systems, Deep Reinforcement Learning, Evolutionary Algorithm"
$ws.Range("G8").Value = "Reinforcement Learning, Multiagent Reinforcement Learning"
$ws.Range("G9").Value = "Model based RL, Diffusion for RL, LLM"
$ws.Range("I9").Value = "Postdoctoral Fellow, Imperial College London"

# Update the sheet view: scroll so column C is the left-most visible column,
# and move the active selection to I10.
$window = $excel.ActiveWindow
$window.ScrollColumn = 3
$ws.Range("I10").Select()
